$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post previously at row 292 ("マヌルネコの珍しい映像...") was removed.
# Deleting the entire row shifts all subsequent rows (293-315) up by one
# (to 292-314) and updates the sheet's used-range dimension accordingly.
$ws.Rows.Item(292).Delete()
